$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "房屋标签" (room tag) header in G3, matching the style of the
# existing header cells in that row (F3).
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "房屋标签"

# Add the empty data cell below it (G4), matching the style used by the
# neighbouring data cells in that row (F4).
$ws.Range("F4").Copy()
$ws.Range("G4").PasteSpecial(-4122)

# Match the column width used in the template for the new column.
$ws.Columns.Item(7).ColumnWidth = 14.140625

# Enable iterative calculation on the workbook.
$wb.IterativeCalculation = $true

# Move/restore the active selection to the new last cell.
$ws.Range("G4").Select()
